# Automatische test-sync: 2025-06-19 21:42:50
$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

$newRow = 25

# Append the new mail-log entry
$logs.Cells.Item($newRow, 1).Value = "Is product X op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik ben geïnteresseerd in product X. Is dit momenteel op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Productinformatie"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 21:42:14"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the conditional formatting ranges so they cover the new row too
$dFc = $logs.Range("D2:D24").FormatConditions
for ($i = 1; $i -le $dFc.Count(); $i++) {
    $dFc.Item($i).ModifyAppliesToRange($logs.Range("D2:D25"))
}

$gFc = $logs.Range("G2:G24").FormatConditions
for ($i = 1; $i -le $gFc.Count(); $i++) {
    $gFc.Item($i).ModifyAppliesToRange($logs.Range("G2:G25"))
}

# Update the Dashboard aggregate count for "Productinformatie" (row 3, col B)
$dash.Cells.Item(3, 2).Value = 5
